{"js": "// Collapse the \"Application satisfies all requirements and functions\n// correctly?  Screen shots of application running in a browser are\n// included?\" checklist text down to just the first sentence, in every\n// place it occurs in the document body (including inside tables).\nconst fullText =\n  \"Application satisfies all requirements and functions correctly?  \" +\n  \"Screen shots of application running in a browser are included?\";\nconst newText = \"Application satisfies all requirements and functions correctly?\";\n\nconst results = context.document.body.search(fullText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Collapse the \"Application satisfies all requirements and functions\n# correctly?  Screen shots of application running in a browser are\n# included?\" checklist text down to just the first sentence, everywhere\n# it appears in the document (it shows up three times, once per lab\n# checklist table).\n$d = $word.ActiveDocument\n\n$oldText = \"Application satisfies all requirements and functions correctly?  Screen shots of application running in a browser are included?\"\n$newText = \"Application satisfies all requirements and functions correctly?\"\n\n$find = $d.Content.Find\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n$find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n"}
